$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-format style from the last existing data row (A18) into the
# brand new row 19 so the new date cell matches the style of the rest of
# column A (custom date number format, centered).
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A2").Value2 = 39400
$ws.Range("B2").Value2 = 2007
$ws.Range("C2").Value2 = 3.14593994906931
$ws.Range("D2").Value2 = 2008
$ws.Range("A3").Value2 = 39765
$ws.Range("B3").Value2 = 2008
$ws.Range("C3").Value2 = 1.769627576887389
$ws.Range("D3").Value2 = 2009
$ws.Range("A4").Value2 = 40130
$ws.Range("B4").Value2 = 2009
$ws.Range("C4").Value2 = -4.774178217057756
$ws.Range("D4").Value2 = 2010
$ws.Range("A5").Value2 = 40494
$ws.Range("B5").Value2 = 2010
$ws.Range("C5").Value2 = 1.97975191822708
$ws.Range("D5").Value2 = 2011
$ws.Range("E5").Value2 = 2.743551941645217
$ws.Range("A6").Value2 = 40862
$ws.Range("B6").Value2 = 2011
$ws.Range("C6").Value2 = 3.452886745653183
$ws.Range("D6").Value2 = 2012
$ws.Range("E6").Value2 = 1.794132456841213
$ws.Range("A7").Value2 = 41228
$ws.Range("B7").Value2 = 2012
$ws.Range("C7").Value2 = 1.239479831392853
$ws.Range("D7").Value2 = 2013
$ws.Range("E7").Value2 = 1.158731032337301
$ws.Range("A8").Value2 = 41592
$ws.Range("B8").Value2 = 2013
$ws.Range("C8").Value2 = 0.2379616621361214
$ws.Range("D8").Value2 = 2014
$ws.Range("E8").Value2 = 1.348985046565354
$ws.Range("A9").Value2 = 41957
$ws.Range("B9").Value2 = 2014
$ws.Range("C9").Value2 = 1.51977456621637
$ws.Range("D9").Value2 = 2015
$ws.Range("E9").Value2 = 0.7739869831243862
$ws.Range("A10").Value2 = 42321
$ws.Range("B10").Value2 = 2015
$ws.Range("C10").Value2 = 1.470039379455756
$ws.Range("D10").Value2 = 2016
$ws.Range("E10").Value2 = 1.319057785023592
$ws.Range("A11").Value2 = 42689
$ws.Range("B11").Value2 = 2016
$ws.Range("C11").Value2 = 1.638797242243251
$ws.Range("D11").Value2 = 2017
$ws.Range("E11").Value2 = 1.236938064849924
$ws.Range("A12").Value2 = 43053
$ws.Range("B12").Value2 = 2017
$ws.Range("C12").Value2 = 2.161565493242668
$ws.Range("D12").Value2 = 2018
$ws.Range("E12").Value2 = 2.320541194291881
$ws.Range("A13").Value2 = 43418
$ws.Range("B13").Value2 = 2018
$ws.Range("C13").Value2 = 2.214251681313772
$ws.Range("D13").Value2 = 2019
$ws.Range("E13").Value2 = 0.6374750548026054
$ws.Range("A14").Value2 = 43783
$ws.Range("B14").Value2 = 2019
$ws.Range("C14").Value2 = 0.6066442151010376
$ws.Range("D14").Value2 = 2020
$ws.Range("E14").Value2 = 0.57214245765278
$ws.Range("A15").Value2 = 44159
$ws.Range("B15").Value2 = 2020
$ws.Range("C15").Value2 = -4.207901339433196
$ws.Range("D15").Value2 = 2021
$ws.Range("E15").Value2 = -0.2586890779524231
$ws.Range("A16").Value2 = 44525
$ws.Range("B16").Value2 = 2021
$ws.Range("C16").Value2 = 1.099928004397532
$ws.Range("D16").Value2 = 2022
$ws.Range("E16").Value2 = 1.794400784768979
$ws.Range("A17").Value2 = 44890
$ws.Range("B17").Value2 = 2022
$ws.Range("C17").Value2 = 2.310042359896225
$ws.Range("D17").Value2 = 2023
$ws.Range("E17").Value2 = 0.8232644777432796
$ws.Range("A18").Value2 = 45254
$ws.Range("B18").Value2 = 2023
$ws.Range("C18").Value2 = 0.0464415346324687
$ws.Range("D18").Value2 = 2024
$ws.Range("E18").Value2 = 0.4457784880425031
$ws.Range("A19").Value2 = 45618
$ws.Range("B19").Value2 = 2024
$ws.Range("C19").Value2 = -0.3101476031197148
$ws.Range("D19").Value2 = 2025
$ws.Range("E19").Value2 = 0.5215192790195111
